$wb = $excel.ActiveWorkbook

# --- Sheet "Servicio": change D4 value from "3 - CENTIMETRO" to "0 - NO APLICA" ---
$wsServicio = $wb.Worksheets.Item("Servicio")
$wsServicio.Range("D4").Value = "0 - NO APLICA"

# --- Update selection on each sheet, and restore "Bien" as the active tab ---
$wsBien = $wb.Worksheets.Item("Bien")
$wsBien.Activate()
$wsBien.Range("A3").Select()

$wsServicio.Activate()
$wsServicio.Range("A4").Select()

$wsBien.Activate()
